$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row: D (Fecha), J (Volumen), K (Precio minimo), L (Precio maximo), M (Precio promedio ponderado), P (Precio $/Kg)
$data = @{
    2 = @{ D = 44505; J = 400; K = 16000; L = 17000; M = 16500; P = 1269 }
    3 = @{ D = 44309; J = 400; K = 26000; L = 27000; M = 26500; P = 2038 }
    4 = @{ D = 44344; J = 400; K = 18500; L = 19000; M = 18750; P = 1442 }
    5 = @{ D = 44414; J = 500; K = 14000; L = 15000; M = 14500; P = 1115 }
    6 = @{ D = 44599; J = 400; K = 15000; L = 16000; M = 15500; P = 1192 }
    7 = @{ D = 44326; J = 460; K = 25000; L = 26000; M = 25500; P = 1962 }
    8 = @{ D = 44426; J = 460; K = 14000; L = 15000; M = 14500; P = 1115 }
    9 = @{ D = 44400; J = 600; K = 15000; L = 16000; M = 15500; P = 1192 }
    10 = @{ D = 44249; J = 400; K = 42000; L = 43000; M = 42500; P = 3269 }
    11 = @{ D = 44260; J = 400; K = 37000; L = 38000; M = 37500; P = 2885 }
    12 = @{ D = 44410; J = 600; K = 14000; L = 15000; M = 14500; P = 1115 }
    13 = @{ D = 44498; J = 400; K = 14000; L = 15000; M = 14500; P = 1115 }
    14 = @{ D = 44442; J = 460; K = 14000; L = 15000; M = 14500; P = 1115 }
    15 = @{ D = 44484; J = 360; K = 14000; L = 15000; M = 14500; P = 1115 }
    16 = @{ D = 44435; J = 480; K = 13000; L = 14000; M = 13500; P = 1038 }
    17 = @{ D = 44333; J = 440; K = 24000; L = 25000; M = 24500; P = 1885 }
    18 = @{ D = 44631; J = 400; K = 16000; L = 17000; M = 16500; P = 1269 }
    19 = @{ D = 44418; J = 500; K = 14000; L = 15000; M = 14500; P = 1115 }
    20 = @{ D = 44365; J = 500; K = 19500; L = 20000; M = 19750; P = 1519 }
    21 = @{ D = 44596; J = 500; K = 16000; L = 17000; M = 16500; P = 1269 }
    22 = @{ D = 44379; J = 600; K = 17000; L = 18000; M = 17500; P = 1346 }
    23 = @{ D = 44644; J = 400; K = 15000; L = 16000; M = 15500; P = 1192 }
    24 = @{ D = 44575; J = 500; K = 14000; L = 15000; M = 14500; P = 1115 }
    25 = @{ D = 44335; J = 480; K = 24500; L = 25000; M = 24750; P = 1904 }
    26 = @{ D = 44419; J = 600; K = 14000; L = 15000; M = 14500; P = 1115 }
    27 = @{ D = 44323; J = 460; K = 25000; L = 26000; M = 25500; P = 1962 }
    28 = @{ D = 44428; J = 480; K = 14000; L = 15000; M = 14500; P = 1115 }
    29 = @{ D = 44412; J = 600; K = 14000; L = 15000; M = 14500; P = 1115 }
    30 = @{ D = 44582; J = 520; K = 15000; L = 16000; M = 15500; P = 1192 }
    31 = @{ D = 44445; J = 600; K = 13000; L = 14000; M = 13500; P = 1038 }
    32 = @{ D = 44383; J = 200; K = 17000; L = 18000; M = 17500; P = 1346 }
    33 = @{ D = 44533; J = 520; K = 17000; L = 18000; M = 17500; P = 1346 }
    34 = @{ D = 44312; J = 400; K = 26000; L = 27000; M = 26500; P = 2038 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("J$row").Value = $vals.J
    $ws.Range("K$row").Value = $vals.K
    $ws.Range("L$row").Value = $vals.L
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("P$row").Value = $vals.P
}
